$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New poll rows (ifop online poll, 11/28), appended as rows 132-134
$newRows = @(
    @{ A=40; B=2021; C=14; D=11; E=24; F="ifop"; G="online"; H="included"; I=1351; J=1; K="T_0.5"; L=7.5; M=2;   N=2; O=7;   P=6; Q=25;               T=10; U=0.5; V=3;   W=19;   X=14 },
    @{ A=40; B=2021; C=14; D=11; E=24; F="ifop"; G="online"; H="included"; I=1351; J=1; K=0.5;     L=8.5; M=2.5; N=2; O=6.5; P=6; Q=25;         S=10;       U=0.5; V=3.5; W=19;   X=15 },
    @{ A=40; B=2021; C=14; D=11; E=24; F="ifop"; G="online"; H="included"; I=1351; J=1; K=0.5;     L=8.5; M=2;   N=2; O=6.5; P=6; Q=25; R=10;               U="T_0.5"; V=4;   W=19.5; X=15 }
)

$startRow = 132
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}

# Scroll / selection bookkeeping to match the new bottom of the sheet
$ws.Activate()
$ws.Range("AB132").Select()
$excel.ActiveWindow.ScrollRow = 124
